$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.39505822365225
$ws.Range("C2").Value = 8.525222701492561
$ws.Range("D2").Value = 5.390011855418382
$ws.Range("E2").Value = 11.57574284004627
$ws.Range("F2").Value = 49.57651527291281
$ws.Range("J2").Value = 10.29394382703299
$ws.Range("K2").Value = 15.97252694999251
$ws.Range("M2").Value = 18.23634088441713
$ws.Range("N2").Value = 24.5982423150714
$ws.Range("B3").Value = 16.23565866884608
$ws.Range("C3").Value = 8.425275859132224
$ws.Range("D3").Value = 5.396439464675019
$ws.Range("E3").Value = 11.58537427125237
$ws.Range("F3").Value = 49.48887899350202
$ws.Range("J3").Value = 10.31127922875067
$ws.Range("K3").Value = 15.87701324052386
$ws.Range("M3").Value = 18.21572541783835
$ws.Range("N3").Value = 24.63127097962127
$ws.Range("B4").Value = 16.14159039517517
$ws.Range("C4").Value = 8.365915977353689
$ws.Range("D4").Value = 5.40131855331923
$ws.Range("E4").Value = 11.59314335261705
$ws.Range("F4").Value = 49.44419366543298
$ws.Range("J4").Value = 10.32301801257107
$ws.Range("K4").Value = 15.82216627869774
$ws.Range("M4").Value = 18.20674977420801
$ws.Range("N4").Value = 24.65339467434958
$ws.Range("B5").Value = 16.10425483877729
$ws.Range("C5").Value = 8.342260119897347
$ws.Range("D5").Value = 5.403541818846189
$ws.Range("E5").Value = 11.59677618763764
$ws.Range("F5").Value = 49.42828629178091
$ws.Range("J5").Value = 10.3280771910496
$ws.Range("K5").Value = 15.80079021640762
$ws.Range("M5").Value = 18.20402125706597
$ws.Range("N5").Value = 24.66287388062727
$ws.Range("B6").Value = 16.0981167495188
$ws.Range("C6").Value = 8.338365141916627
$ws.Range("D6").Value = 5.403925197621162
$ws.Range("E6").Value = 11.59740762248292
$ws.Range("F6").Value = 49.42578413543126
$ws.Range("J6").Value = 10.32893391189463
$ws.Range("K6").Value = 15.79730012768233
$ws.Range("M6").Value = 18.20362438668594
$ws.Range("N6").Value = 24.66447589549005
$ws.Range("B7").Value = 16.1410827805129
$ws.Range("C7").Value = 8.365594747687268
$ws.Range("D7").Value = 5.401347584913045
$ws.Range("E7").Value = 11.59319045559868
$ws.Range("F7").Value = 49.44396980210765
$ws.Range("J7").Value = 10.32308512650862
$ws.Range("K7").Value = 15.82187402361793
$ws.Range("M7").Value = 18.20670921074456
$ws.Range("N7").Value = 24.65352063719082
$ws.Range("B8").Value = 16.33933335829825
$ws.Range("C8").Value = 8.490361969501469
$ws.Range("D8").Value = 5.392034872055737
$ws.Range("E8").Value = 11.57867890912999
$ws.Range("F8").Value = 49.54440998367918
$ws.Range("J8").Value = 10.29969399727489
$ws.Range("K8").Value = 15.93881741590485
$ws.Range("M8").Value = 18.22847048921924
$ws.Range("N8").Value = 24.6092478494735
$ws.Range("B9").Value = 16.75636212660856
$ws.Range("C9").Value = 8.749640890341674
$ws.Range("D9").Value = 5.381147751617591
$ws.Range("E9").Value = 11.56492542711219
$ws.Range("F9").Value = 49.81334306824603
$ws.Range("J9").Value = 10.2625002116527
$ws.Range("K9").Value = 16.19736501944233
$ws.Range("M9").Value = 18.30019734324537
$ws.Range("N9").Value = 24.53706314533298
$ws.Range("B10").Value = 17.07729917188323
$ws.Range("C10").Value = 8.947217428870426
$ws.Range("D10").Value = 5.377610476326076
$ws.Range("E10").Value = 11.56375586316005
$ws.Range("F10").Value = 50.05413698946136
$ws.Range("J10").Value = 10.24044919917012
$ws.Range("K10").Value = 16.40380482387852
$ws.Range("M10").Value = 18.37034949442822
$ws.Range("N10").Value = 24.49295535631567
$ws.Range("B11").Value = 17.22586761958407
$ws.Range("C11").Value = 9.038242865189755
$ws.Range("D11").Value = 5.37696164040238
$ws.Range("E11").Value = 11.56515551078176
$ws.Range("F11").Value = 50.17288486241435
$ws.Range("J11").Value = 10.23156002572408
$ws.Range("K11").Value = 16.50099017372074
$ws.Range("M11").Value = 18.40598589497818
$ws.Range("N11").Value = 24.47482904465743
$ws.Range("B12").Value = 17.28244610970782
$ws.Range("C12").Value = 9.072843642523697
$ws.Range("D12").Value = 5.376853225500351
$ws.Range("E12").Value = 11.56596233606286
$ws.Range("F12").Value = 50.21915784560939
$ws.Range("J12").Value = 10.22835787249782
$ws.Range("K12").Value = 16.53823512940017
$ws.Range("M12").Value = 18.42000917043437
$ws.Range("N12").Value = 24.46824396038389
$ws.Range("B13").Value = 17.2702475527388
$ws.Range("C13").Value = 9.065386430792497
$ws.Range("D13").Value = 5.376870482185145
$ws.Range("E13").Value = 11.56577627740261
$ws.Range("F13").Value = 50.20913437115427
$ws.Range("J13").Value = 10.22904022396392
$ws.Range("K13").Value = 16.53019451391102
$ws.Range("M13").Value = 18.41696561697125
$ws.Range("N13").Value = 24.46964976628758
$ws.Range("B14").Value = 17.23051626872378
$ws.Range("C14").Value = 9.041087048783774
$ws.Range("D14").Value = 5.376949974443404
$ws.Range("E14").Value = 11.56521634710816
$ws.Range("F14").Value = 50.1766657119726
$ws.Range("J14").Value = 10.2312932974993
$ws.Range("K14").Value = 16.50404562373213
$ws.Range("M14").Value = 18.40712904739027
$ws.Range("N14").Value = 24.47428169398224
$ws.Range("B15").Value = 17.20621971554718
$ws.Range("C15").Value = 9.026219101866683
$ws.Range("D15").Value = 5.37701651953432
$ws.Range("E15").Value = 11.56490939148091
$ws.Range("F15").Value = 50.1569471866949
$ws.Range("J15").Value = 10.232694718761
$ws.Range("K15").Value = 16.48808551836724
$ws.Range("M15").Value = 18.40117248018343
$ws.Range("N15").Value = 24.47715521764895
$ws.Range("B16").Value = 17.06763727638203
$ws.Range("C16").Value = 8.941288910200464
$ws.Range("D16").Value = 5.377672135372411
$ws.Range("E16").Value = 11.56370318385756
$ws.Range("F16").Value = 50.04656030631167
$ws.Range("J16").Value = 10.24105308724868
$ws.Range("K16").Value = 16.39751706114895
$ws.Range("M16").Value = 18.36809492574407
$ws.Range("N16").Value = 24.49417897862668
$ws.Range("B17").Value = 16.98324403330273
$ws.Range("C17").Value = 8.88945681017446
$ws.Range("D17").Value = 5.378319733818808
$ws.Range("E17").Value = 11.56345732070496
$ws.Range("F17").Value = 49.9811876296212
$ws.Range("J17").Value = 10.24647299653719
$ws.Range("K17").Value = 16.34277540055666
$ws.Range("M17").Value = 18.34875193894442
$ws.Range("N17").Value = 24.5051191514309
$ws.Range("B18").Value = 16.93494849097583
$ws.Range("C18").Value = 8.859754300381685
$ws.Range("D18").Value = 5.378782653834124
$ws.Range("E18").Value = 11.5634977923832
$ws.Range("F18").Value = 49.94445505880805
$ws.Range("J18").Value = 10.24969788419579
$ws.Range("K18").Value = 16.31159942055047
$ws.Range("M18").Value = 18.33797732429346
$ws.Range("N18").Value = 24.51159410257001
$ws.Range("B19").Value = 16.91864012132312
$ws.Range("C19").Value = 8.849717449153061
$ws.Range("D19").Value = 5.378954947042629
$ws.Range("E19").Value = 11.56354276170218
$ws.Range("F19").Value = 49.93216767038049
$ws.Range("J19").Value = 10.25080824556202
$ws.Range("K19").Value = 16.30109784558216
$ws.Range("M19").Value = 18.33438970392435
$ws.Range("N19").Value = 24.51381774025199
$ws.Range("B20").Value = 16.99220282872944
$ws.Range("C20").Value = 8.894963290668453
$ws.Range("D20").Value = 5.378241441451371
$ws.Range("E20").Value = 11.56346467537581
$ws.Range("F20").Value = 49.98805695083229
$ws.Range("J20").Value = 10.24588491381205
$ws.Range("K20").Value = 16.34857086239095
$ws.Range("M20").Value = 18.35077475829693
$ws.Range("N20").Value = 24.50393566684655
$ws.Range("B21").Value = 17.24217806401092
$ws.Range("C21").Value = 9.048221064890543
$ws.Range("D21").Value = 5.376922906440705
$ws.Range("E21").Value = 11.56537330802837
$ws.Range("F21").Value = 50.18616726315032
$ws.Range("J21").Value = 10.23062706656099
$ws.Range("K21").Value = 16.51171439117716
$ws.Range("M21").Value = 18.4100040001025
$ws.Range("N21").Value = 24.47291361216674
$ws.Range("B22").Value = 17.40738613245049
$ws.Range("C22").Value = 9.149134643911344
$ws.Range("D22").Value = 5.376860958288924
$ws.Range("E22").Value = 11.5682336252155
$ws.Range("F22").Value = 50.32324487294807
$ws.Range("J22").Value = 10.22161089307628
$ws.Range("K22").Value = 16.62090660858221
$ws.Range("M22").Value = 18.45179084824575
$ws.Range("N22").Value = 24.45426496992797
$ws.Range("B23").Value = 17.31906046751383
$ws.Range("C23").Value = 9.095217441733709
$ws.Range("D23").Value = 5.376821120165078
$ws.Range("E23").Value = 11.56655979081362
$ws.Range("F23").Value = 50.24939504236414
$ws.Range("J23").Value = 10.2263356240086
$ws.Range("K23").Value = 16.56240308225706
$ws.Range("M23").Value = 18.42920931345714
$ws.Range("N23").Value = 24.46406925136272
$ws.Range("B24").Value = 16.98815185801108
$ws.Range("C24").Value = 8.892473507846791
$ws.Range("D24").Value = 5.378276555151912
$ws.Range("E24").Value = 11.56346078391761
$ws.Range("F24").Value = 49.98494867834281
$ws.Range("J24").Value = 10.24615044675661
$ws.Range("K24").Value = 16.34594981081837
$ws.Range("M24").Value = 18.34985916363459
$ws.Range("N24").Value = 24.50447014301685
$ws.Range("B25").Value = 16.64079141866217
$ws.Range("C25").Value = 8.678121090376061
$ws.Range("D25").Value = 5.383306808568877
$ws.Range("E25").Value = 11.56707445192781
$ws.Range("F25").Value = 49.73294633787103
$ws.Range("J25").Value = 10.27163468487228
$ws.Range("K25").Value = 16.12442470455842
$ws.Range("M25").Value = 18.27770701164925
$ws.Range("N25").Value = 24.55502383133999
